# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" (fund-level detail, columns A-H)
#    right after "2021-Q4" and right before "总计".
# 2) Insert a new summary row for "2022-Q1" at the top of the "总计" sheet's
#    data block, re-indexing column A for the rows that shift down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Part 1: new "2022-Q1" worksheet
# ---------------------------------------------------------------------------

$template = $wb.Worksheets.Item("2021-Q4")
$anchor = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $anchor)
$newSheet.Name = "2022-Q1"

# Clone the formatting (header style + bordered/bold index column) from the
# same-shaped "2021-Q4" sheet so the new sheet's styles line up (s="2" on
# the header row and on column A).
$template.Range("A1:H11").Copy()
$newSheet.Range("A1:H11").PasteSpecial(-4122)
$template.Range("A11:H11").Copy()
$newSheet.Range("A12:H12").PasteSpecial(-4122)
$newSheet.Range("A13:H13").PasteSpecial(-4122)
$newSheet.Range("A14:H14").PasteSpecial(-4122)
$newSheet.Range("A15:H15").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$fundRows = @(
    @(0, "003713", "英大睿盛灵活配置混合A", "5.99", "87.42", "7.01", "0.4199", 3),
    @(1, "003714", "英大睿盛灵活配置混合C", "2.40", "87.42", "7.01", "0.1682", 3),
    @(2, "002067", "诺安精选回报灵活配置混合", "4.89", "25.39", "1.23", "0.0601", 3),
    @(3, "002305", "光大保德信风格轮动混合A", "3.99", "90.94", "1.45", "0.0579", 9),
    @(4, "000714", "诺安稳健回报灵活配置混合A", "1.96", "64.55", "2.94", "0.0576", 2),
    @(5, "004634", "新疆前海联合泳涛灵活配置混合A", "1.33", "89.65", "4.15", "0.0552", 9),
    @(6, "002052", "诺安稳健回报灵活配置混合C", "1.63", "64.55", "2.94", "0.0479", 2),
    @(7, "007499", "光大保德信风格轮动混合C", "2.82", "90.94", "1.45", "0.0409", 9),
    @(8, "003446", "英大睿鑫灵活配置混合A", "0.59", "89.46", "6.16", "0.0363", 4),
    @(9, "003447", "英大睿鑫灵活配置混合C", "0.51", "89.46", "6.16", "0.0314", 4),
    @(10, "002145", "诺安景鑫灵活配置混合", "0.53", "83.45", "4.01", "0.0213", 6),
    @(11, "320016", "诺安多策略混合", "0.19", "80.02", "4.25", "0.0081", 4),
    @(12, "001608", "英大策略优选混合C", "0.03", "89.86", "7.59", "0.0023", 3),
    @(13, "007041", "新疆前海联合泳涛灵活配置混合C", "0.00", "89.65", "4.15", "0", 9)
)

$r = 2
foreach ($row in $fundRows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]

    $newSheet.Cells.Item($r, 2).Value = "'" + $row[1]
    $newSheet.Cells.Item($r, 2).Style = "Normal"

    $newSheet.Cells.Item($r, 3).Value = $row[2]

    $newSheet.Cells.Item($r, 4).Value = "'" + $row[3]
    $newSheet.Cells.Item($r, 4).Style = "Normal"

    $newSheet.Cells.Item($r, 5).Value = "'" + $row[4]
    $newSheet.Cells.Item($r, 5).Style = "Normal"

    $newSheet.Cells.Item($r, 6).Value = "'" + $row[5]
    $newSheet.Cells.Item($r, 6).Style = "Normal"

    if ($r -eq 15) {
        # last row's held-value column is stored as a real number (0), not text
        $newSheet.Cells.Item($r, 7).Value = 0
    } else {
        $newSheet.Cells.Item($r, 7).Value = "'" + $row[6]
        $newSheet.Cells.Item($r, 7).Style = "Normal"
    }

    $newSheet.Cells.Item($r, 8).Value = $row[7]

    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Part 2: new summary row in "总计"
# ---------------------------------------------------------------------------

$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 14
$total.Range("D2").Value = 1.01

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
